# daily auto push: 2026-02-22 09:43 UTC
# Insert a new data row for 2026/02/22 (日, hour 16) right before the
# existing row 850 (2026/12/29) entry, shifting every subsequent row
# down by one. This extends the used range from A1:D891 to A1:D892.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 850 (and everything below it) down by one row.
$ws.Rows.Item(850).Insert()

# The new row 850 holds the inserted reading. Force column A to stay a
# plain text date-like label (matching every other row in the sheet)
# instead of letting Excel auto-convert the "yyyy/mm/dd" string into a
# real date serial number.
$ws.Range("A850").NumberFormat = "@"
$ws.Range("A850").Value = "2026/02/22"
$ws.Range("A850").ClearFormats()

$ws.Range("B850").Value = "日"
$ws.Range("C850").Value = 16
$ws.Range("D850").Value = 201
